$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.930.45'
$ws.Range('E2').Value = '  +6.18%  '
$ws.Range('D3').Value = '3.480.02'
$ws.Range('E3').Value = '  +4.14%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '412.01'
$ws.Range('E5').Value = '  +0.29%  '
$ws.Range('D6').Value = '129.78'
$ws.Range('E6').Value = '  +14.20%  '
$ws.Range('D7').Value = '3.471.02'
$ws.Range('E7').Value = '  +4.10%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.600'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.89%  '
$ws.Range('E9').Value = '  +0.12%  '
$ws.Range('D10').Value = '0.699'
$ws.Range('E10').Value = '  +8.79%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.130'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +30.73%  '
$ws.Range('D12').Value = '43.39'
$ws.Range('E12').Value = '  +7.61%  '
$ws.Range('E13').Value = '  -0.53%  '
$ws.Range('D14').Value = '4.025.27'
$ws.Range('E14').Value = '  +3.87%  '
$ws.Range('D15').Value = '8.76'
$ws.Range('E15').Value = '  +2.45%  '
$ws.Range('D16').Value = '20.22'
$ws.Range('E16').Value = '  +3.56%  '
$ws.Range('D17').Value = '3.429.96'
$ws.Range('E17').Value = '  +2.11%  '
$ws.Range('D18').Value = '62.846.23'
$ws.Range('E18').Value = '  +6.31%  '
$ws.Range('E19').Value = '  -0.46%  '
$ws.Range('E20').Value = '  +2.58%  '
$ws.Range('E21').Value = '  +23.54%  '
$ws.Range('E22').Value = '  -0.80%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '82.60'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +8.98%  '
$ws.Range('E24').Value = '  -0.14%  '
$ws.Range('D25').Value = '314.47'
$ws.Range('E25').Value = '  +3.10%  '
$ws.Range('E26').Value = '  -0.92%  '
$ws.Range('D27').Value = '30.52'
$ws.Range('E27').Value = '  +6.25%  '
$ws.Range('D28').Value = '8.23'
$ws.Range('E28').Value = '  +3.69%  '
$ws.Range('E29').Value = '  +1.99%  '
$ws.Range('E30').Value = '  +2.54%  '
$ws.Range('D31').Value = '4.36'
$ws.Range('E32').Value = '  +1.44%  '
$ws.Range('D33').Value = '12.16'
$ws.Range('E33').Value = '  +5.65%  '
$ws.Range('D34').Value = '44.44'
$ws.Range('E34').Value = '  +10.71%  '
$ws.Range('E35').Value = '  +25.31%  '
$ws.Range('E36').Value = '  +0.09%  '
$ws.Range('D37').Value = '0.0494'
$ws.Range('E37').Value = '  -5.45%  '
$ws.Range('E38').Value = '  +1.04%  '
$ws.Range('E39').Value = '  +5.25%  '
$ws.Range('E40').Value = '  -0.42%  '
$ws.Range('D41').Value = '3.05'
$ws.Range('E41').Value = '  -3.02%  '
$ws.Range('D42').Value = '2.01'
$ws.Range('E42').Value = '  +3.38%  '
$ws.Range('E43').Value = '  +2.46%  '
$ws.Range('D44').Value = '138.18'
$ws.Range('E44').Value = '  +0.39%  '
$ws.Range('D45').Value = '17.84'
$ws.Range('E45').Value = '  +4.07%  '
$ws.Range('E46').Value = '  +2.32%  '
$ws.Range('D47').Value = '3.99'
$ws.Range('E47').Value = '  +0.01%  '
$ws.Range('E48').Value = '  +0.77%  '
$ws.Range('D49').Value = '22.52'
$ws.Range('E49').Value = '  -0.05%  '
$ws.Range('D50').Value = '2.218.07'
$ws.Range('E50').Value = '  +0.13%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '3.822.15'
$ws.Range('E51').Value = '  +3.91%  '
